$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UK")
$ws.Activate()

# Add new regressor row below the existing data (S3e - social care provision,
# used to align with social care receipt)
$ws.Range("A25").Value = "S3e"
$ws.Range("B25").Value = 1.2428999999999999

# Update view state: scroll the window down and move the active selection
# to the next empty row below the newly added data
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$ws.Range("B26").Select()
